# CardData2ndDraft.xlsx edit: add a new "Method Overloading" category column,
# switch Yes/No style headers to numeric-count headers, and append ten new
# programming-language rows (11-20) with their category counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1) - relabel columns B:G and add a new H column.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Error Handling"
$ws.Range("C1").Value = "Garbage Collection"
$ws.Range("D1").Value = "Variables & Classes"
$ws.Range("E1").Value = "Regular Expressions"
$ws.Range("F1").Value = "Language Integration"
$ws.Range("G1").Value = "Built-In Security"

# New column header H1 - copy the header formatting from G1, then set value.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Method Overloading"

# ---------------------------------------------------------------------
# 2. Existing data rows 2-10: updated numeric values.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 3
$ws.Range("E2").Value = 0

$ws.Range("B3").Value = 2
$ws.Range("E3").Value = 1

$ws.Range("E4").Value = 2

$ws.Range("E5").Value = 4

$ws.Range("E6").Value = 4

$ws.Range("E7").Value = 3

$ws.Range("E8").Value = 3

$ws.Range("E10").Value = 3

# New column H for rows 2-4 (copy number-cell format from column G, then set value).
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Value = 0

$ws.Range("G3").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Value = 0

$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = 0

# New highlighted (green) cell I2.
$ws.Range("I2").Interior.Color = 5296274

# ---------------------------------------------------------------------
# 3. New rows 11-20: more programming languages + category counts.
#    Copy formatting from the last existing rows so fonts/fills/number
#    styles match (style ids 1 for labels, 4 for number cells) - only
#    onto the specific cells that end up with content, so no stray
#    empty-but-styled cells are introduced.
# ---------------------------------------------------------------------
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A20").PasteSpecial(-4122) | Out-Null

$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11:C12").PasteSpecial(-4122) | Out-Null

$ws.Range("D10").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null

$ws.Range("E10").Copy() | Out-Null
$ws.Range("E11:E13").PasteSpecial(-4122) | Out-Null

$ws.Range("F10").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null

$ws.Range("G10").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null

# Row 11 - C
$ws.Range("A11").Value = "C"
$ws.Range("C11").Value = 4
$ws.Range("E11").Value = 4

# Row 12 - PHP
$ws.Range("A12").Value = "PHP"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 2

# Row 13 - Java Script
$ws.Range("A13").Value = "Java Script"
$ws.Range("E13").Value = 2

# Row 14 - Scheme
$ws.Range("A14").Value = "Scheme"

# Row 15 - Haskell
$ws.Range("A15").Value = "Haskell"

# Row 16 - Pascal
$ws.Range("A16").Value = "Pascal"

# Row 17 - Scala
$ws.Range("A17").Value = "Scala"

# Row 18 - Go
$ws.Range("A18").Value = "Go"

# Row 19 - Perl
$ws.Range("A19").Value = "Perl"

# Row 20 - Objective C
$ws.Range("A20").Value = "Objective C"

# ---------------------------------------------------------------------
# 4. Column widths / row heights (best effort match of the saved layout).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 10.5
$ws.Columns.Item(8).ColumnWidth = 9.833333333333334

$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(20).RowHeight = 19.5

# ---------------------------------------------------------------------
# 5. Selection, matching the saved workbook's active cell.
# ---------------------------------------------------------------------
$ws.Range("D3").Select() | Out-Null
